$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (Coin name) updates ---
$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("B47").Value = 'RenderToken'

# --- Column C (Link) updates ---
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'

# --- Column E (Volume %) updates: plain text (leading/trailing spaces
# and "%" keep Excel from re-parsing these as numbers) ---
$ws.Range("E2").Value = '  +1.27%  '
$ws.Range("E3").Value = '  +1.47%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("E5").Value = '  -0.22%  '
$ws.Range("E6").Value = '  -0.11%  '
$ws.Range("E7").Value = '  -1.77%  '
$ws.Range("E8").Value = '  -0.53%  '
$ws.Range("E9").Value = '  -0.24%  '
$ws.Range("E10").Value = '  +2.94%  '
$ws.Range("E11").Value = '  -2.17%  '
$ws.Range("E12").Value = '  -1.10%  '
$ws.Range("E13").Value = '  +0.60%  '
$ws.Range("E14").Value = '  -0.86%  '
$ws.Range("E15").Value = '  +0.67%  '
$ws.Range("E16").Value = '  -1.94%  '
$ws.Range("E17").Value = '  -0.09%  '
$ws.Range("E18").Value = '  -0.30%  '
$ws.Range("E19").Value = '  -0.11%  '
$ws.Range("E20").Value = '  -0.08%  '
$ws.Range("E21").Value = '  +1.09%  '
$ws.Range("E22").Value = '  -0.89%  '
$ws.Range("E23").Value = '  -0.81%  '
$ws.Range("E24").Value = '  -0.25%  '
$ws.Range("E25").Value = '  +0.55%  '
$ws.Range("E26").Value = '  -0.44%  '
$ws.Range("E27").Value = '  -0.10%  '
$ws.Range("E28").Value = '  -0.89%  '
$ws.Range("E29").Value = '  -1.49%  '
$ws.Range("E30").Value = '  -1.05%  '
$ws.Range("E31").Value = '  +0.14%  '
$ws.Range("E32").Value = '  -1.97%  '
$ws.Range("E33").Value = '  -0.36%  '
$ws.Range("E34").Value = '  -0.35%  '
$ws.Range("E35").Value = '  -1.14%  '
$ws.Range("E36").Value = '  -1.65%  '
$ws.Range("E37").Value = '  -0.56%  '
$ws.Range("E38").Value = '  -2.11%  '
$ws.Range("E39").Value = '  -0.97%  '
$ws.Range("E40").Value = '  -0.54%  '
$ws.Range("E41").Value = '  +0.54%  '
$ws.Range("E42").Value = '  -2.32%  '
$ws.Range("E43").Value = '  -1.67%  '
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("E45").Value = '  -0.48%  '
$ws.Range("E46").Value = '  -0.12%  '
$ws.Range("E47").Value = '  +0.19%  '
$ws.Range("E48").Value = '  -0.80%  '
$ws.Range("E49").Value = '  -1.76%  '
$ws.Range("E50").Value = '  +4.06%  '
$ws.Range("E51").Value = '  -0.09%  '

# --- Column D (Price) updates ---
# Several of these look like plain numbers (e.g. "1.005"); a bare string
# assignment that parses as a number gets auto-converted by Excel into a
# numeric cell, which would change the stored type. Force text storage per
# cell via NumberFormat "@", then restore the default "Normal" style so no
# stray cell style lingers (original cells carry no explicit style either).
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '28.522.60'
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.867.76'
$c.Style = "Normal"
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.005'
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '325.03'
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '1.005'
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.4558'
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3840'
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.07835'
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.9905'
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '21.58'
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '1.849.69'
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '6.914'
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '5.641'
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.06947'
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '86.79'
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '1.006'
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.000009925'
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '16.68'
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '28.492.72'
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '5.256'
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '10.93'
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '2.098'
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.090.28'
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '153.81'
$c.Style = "Normal"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '5.691'
$c.Style = "Normal"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '1.938'
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '117.51'
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.09276'
$c.Style = "Normal"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '0.9098'
$c.Style = "Normal"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '5.263'
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '1.321'
$c.Style = "Normal"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '3.294'
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.05709'
$c.Style = "Normal"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '1.138'
$c.Style = "Normal"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.02054'
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '7.671'
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.5558'
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '9.654'
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.07101'
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '11.58'
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.5243'
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '1.131'
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '2.121'
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '1.816'
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '111.76'
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '2.417'
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '1.005'
$c.Style = "Normal"
